# Update the "想去人数" (want-to-go count) figures in the F column for the
# two sheets that carry the exhibition listing data: "展览" and "全部类型".
# Row -> new F value mapping (per the diff):
#   F3  -> 8097
#   F8  -> 139
#   F10 -> 191
#   F12 -> 731
#   F13 -> 187
#   F14 -> 2807
#   F20 -> 93

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 8097
    8  = 139
    10 = 191
    12 = 731
    13 = 187
    14 = 2807
    20 = 93
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
